$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# First-page header picture (BTEC logo): image1.jpg -> image2.jpg
$headerShape = $sec.Headers.Item(2).Range.InlineShapes.Item(1).Range.InlineShapes.Item(1)
$headerShape.Name = "image2.jpg"

# First-page footer picture (Pearson logo): image2.png -> image1.png
$footerShape1 = $sec.Footers.Item(1).Range.InlineShapes.Item(1).Range.InlineShapes.Item(1)
$footerShape1.Name = "image1.png"

# Default footer picture (Pearson logo): image2.png -> image1.png
$footerShape2 = $sec.Footers.Item(2).Range.InlineShapes.Item(1).Range.InlineShapes.Item(1)
$footerShape2.Name = "image1.png"
